$wb = $excel.ActiveWorkbook

# --- "cim" worksheet: insert a new "skeleton" row for /courseleaf/fonts ---
$wsCim = $wb.Worksheets.Item("cim")

# Insert a new row at position 6 (after the "/email" skeleton row, before "cgis")
$wsCim.Rows("6:6").Insert()

# Copy formatting for the new row from the row above (A:C) and from the row
# that used to be directly below (D), matching the target layout exactly.
$wsCim.Range("A5:C5").Copy()
$wsCim.Range("A6:C6").PasteSpecial(-4122)
$wsCim.Range("D7").Copy()
$wsCim.Range("D6").PasteSpecial(-4122)

$wsCim.Range("A6").Value = "skeleton"
$wsCim.Range("B6").Value = "/web/courseleaf/fonts"
$wsCim.Range("C6").Value = "/web/<progDir>/fonts"

$wsCim.Range("C24").Select()

# --- "cat" worksheet: remove the now-redundant focussearch skeleton row ---
$wsCat = $wb.Worksheets.Item("cat")
$wsCat.Rows("12:12").Delete()

$wsCat.Range("B24").Select()

# --- "-Instructions" worksheet: leftover cursor position from the edit session ---
$wsInstructions = $wb.Worksheets.Item("-Instructions")
$wsInstructions.Range("G4").Select()
